$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrección columna posición de las clasificaciones
# The "position" values in column A (rows 2-7) were off by one (1-indexed
# instead of 0-indexed); decrement each by 1.
for ($r = 2; $r -le 7; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = $cell.Value2 - 1
}
